# edit.ps1
# Applies the "want-to-go count" (想去人数, column F) increments produced
# by the gh-pages scraper re-run (commit "Update gh-pages to output
# generated at 456a3b4") to the three affected worksheets:
#   展览     (Exhibitions)
#   本地生活 (Local life)
#   全部类型 (All types)
# The 演出 (Performances) sheet is untouched by this commit.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet "展览") ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F7").Value  = 12753   # was 12747
$wsExpo.Range("F10").Value = 2840    # was 2838
$wsExpo.Range("F12").Value = 6094    # was 6089
$wsExpo.Range("F25").Value = 85      # was 84
$wsExpo.Range("F31").Value = 6376    # was 6373
$wsExpo.Range("F33").Value = 157     # was 156
$wsExpo.Range("F41").Value = 211     # was 210
$wsExpo.Range("F44").Value = 113     # was 112
$wsExpo.Range("F46").Value = 1696    # was 1695
$wsExpo.Range("F48").Value = 140     # was 139

# --- 本地生活 (sheet "本地生活") ---
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 395     # was 394

# --- 全部类型 (sheet "全部类型") ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 395      # was 394
$wsAll.Range("F9").Value  = 12753    # was 12747
$wsAll.Range("F12").Value = 2840     # was 2838
$wsAll.Range("F14").Value = 6094     # was 6089
$wsAll.Range("F31").Value = 6376     # was 6373
$wsAll.Range("F34").Value = 157      # was 156
$wsAll.Range("F40").Value = 211      # was 210
$wsAll.Range("F46").Value = 1696     # was 1695
$wsAll.Range("F48").Value = 140      # was 139
